$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Full results")
$ws2 = $wb.Worksheets.Item("For plotting")

# --- Sheet "Full results" ---
$ws1.Range("C2").Value = 0.933934515757714
$ws1.Range("D2").Value = 0.0661353206243123
$ws1.Range("E2").Value = 1.00006983638203

$ws1.Range("J2").Value = 0.0661307022953231
$ws1.Range("K2").Value = 0.0623069777502883
$ws1.Range("L2").Value = 0.0807980086899275
$ws1.Range("M2").Value = 0.0825849676727291
$ws1.Range("N2").Value = 0.143104986440216

$ws1.Range("F3").Value = 0.932147431980163
$ws1.Range("G3").Value = 0.0623113290441894

$ws1.Range("H4").Value = 0.851343780649633
$ws1.Range("I4").Value = 0.026012353861055
$ws1.Range("O4").Value = 0.148715669968052

# --- Sheet "For plotting" ---
$ws2.Range("C2").Value = 0.0661307022953231
$ws2.Range("D2").Value = -0.00190399304253661
$ws2.Range("E2").Value = 0.134165397633183
$ws2.Range("F2").Value = 948

$ws2.Range("C3").Value = 0.143104986440216
$ws2.Range("D3").Value = 0.081975965764306
$ws2.Range("E3").Value = 0.204234007116126
$ws2.Range("F3").Value = 948

$ws2.Range("C4").Value = 0.148715669968052
$ws2.Range("D4").Value = 0.0756297148128868
$ws2.Range("E4").Value = 0.221801625123218
$ws2.Range("F4").Value = 948
